# Add a new "blocked_domains" worksheet right after the existing "sources"
# sheet, populate it with a small list of blocked domains, and leave it as
# the active sheet/tab (matching the author's in-app edit).

$wb = $excel.ActiveWorkbook
$sourcesSheet = $wb.Worksheets.Item(1)

# Insert the new sheet immediately after "sources" so tab order becomes
# sources, blocked_domains.
$ws = $wb.Worksheets.Add($null, $sourcesSheet)
$ws.Name = "blocked_domains"

$ws.Range("A1").Value = "Domains"
$ws.Range("A2").Value = "finance.yahoo.com"
$ws.Range("A3").Value = "bloomberg.com"

# Mirror Excel leaving the cursor one row below the last entry after typing
# values in and pressing Enter, and keep the new sheet active/selected.
[void]$ws.Range("A4").Select()
$ws.Activate()
